# Add two new columns "I0" (I) and "IF" (J) to the worksheet, mirroring the
# existing header/body formatting used by column H ("IP").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (bold header style incl. borders for row 1, and the
# plain body style for rows 2:17) from column H into the two new columns so
# that no new style entries are introduced.
$ws.Range("H1:H17").Copy()
$ws.Range("I1:I17").PasteSpecial(-4122)

$ws.Range("H1:H17").Copy()
$ws.Range("J1:J17").PasteSpecial(-4122)

# Header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Body values.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 11
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 12
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 5
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 6
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 8
$ws.Range("I8").Value = 9
$ws.Range("J8").Value = 9
$ws.Range("I9").Value = 9
$ws.Range("J9").Value = 9
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("I11").Value = 6
$ws.Range("J11").Value = 6
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 3
$ws.Range("I13").Value = 9
$ws.Range("J13").Value = 9
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 5
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 8
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 4
